$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- New BOM line (row 14): TCA9517 i2c buffer header-breakaway flavour ---
# Values first (string insertion order matters for shared-string indices:
# Fournisseur ref, Repere, Ref. Fabricant, Description).
$ws.Range("G14").Value = "538-22-28-4205"
$ws.Range("A14").Value = "Headers 20P - to break away"
$ws.Range("C14").Value = "22-28-4205"
$ws.Range("D14").Value = "Header breakaway 20P 1R"
$ws.Range("E14").Value = "Molex"
$ws.Range("F14").Value = "Mouser"
$ws.Range("B14").Value = 2
$ws.Range("H14").Value = 0.918
$ws.Range("I14").Formula = "=H14*B14"

# Formats: reuse existing styles from neighbouring rows so no duplicate
# style/numFmt entries are minted.
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial($xlPasteFormats)
$ws.Range("D13").Copy()
$ws.Range("D14").PasteSpecial($xlPasteFormats)
$ws.Range("F13").Copy()
$ws.Range("F14").PasteSpecial($xlPasteFormats)
$ws.Range("H13").Copy()
$ws.Range("H14").PasteSpecial($xlPasteFormats)
$ws.Range("I13").Copy()
$ws.Range("I14").PasteSpecial($xlPasteFormats)

# A14 gets a brand-new italic style (no existing donor cell for it).
$ws.Range("A14").Font.Italic = $true

# Remember format donors from the old summary rows (15-18) before removing
# them -- copy each into a scratch column far out of the way first so the
# upcoming row Delete (which shifts everything below up) cannot clobber the
# destinations before we use them.
$ws.Range("I15").Copy()
$ws.Range("Z1").PasteSpecial($xlPasteFormats)
$ws.Range("I16").Copy()
$ws.Range("Z2").PasteSpecial($xlPasteFormats)
$ws.Range("I17").Copy()
$ws.Range("Z3").PasteSpecial($xlPasteFormats)
$ws.Range("J18").Copy()
$ws.Range("Z4").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

# Remove the now-stale old summary rows (15-18) -- nothing of ours lives
# below row 18 yet, so this cannot disturb any new content.
$ws.Range("A15:J18").Delete()

# --- Rebuild the summary block further down the sheet (rows 20-23) ---
$ws.Range("Z1").Copy()
$ws.Range("I20").PasteSpecial($xlPasteFormats)
$ws.Range("Z2").Copy()
$ws.Range("I21").PasteSpecial($xlPasteFormats)
$ws.Range("Z3").Copy()
$ws.Range("I22").PasteSpecial($xlPasteFormats)
$ws.Range("Z4").Copy()
$ws.Range("J23").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

$ws.Range("H20").Value = "Mouser"
$ws.Range("I20").Formula = "=SUM(I4:I14,I2)"

$ws.Range("H21").Value = "Digikey"
$ws.Range("I21").Formula = "=I3"

$ws.Range("H22").Value = "Total"
$ws.Range("I22").Formula = "=SUM(I2:I14)"

# Clean up the scratch cells used to stage formats.
$ws.Range("Z1:Z4").Delete()

$ws.Range("C16").Select()
